# "adding new profile test cases"
#
# The "D suite.xlsx" / "Test Cases" sheet gets two new profile test-case
# results swapped:
#   - E7  (ProfileUpdateTest / OPQA-495)                  FAIL -> SKIP
#   - E12 (ProfileInterestSkillsUpdateTest / OPQA-496)     SKIP -> PASS
# (the shared "FAIL" string itself is renamed to "PASS" and re-used by E12,
# while E7 now points at the existing "SKIP" string.)
#
# The sheet's view also scrolled/re-selected from C15 (with A8 at the top)
# down to D20 (with A9 at the top).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- data edits -----------------------------------------------------------
# Row 7 held the "FAIL" result; row 12 held "SKIP". Swap them so row 7
# becomes "SKIP" and row 12 becomes "PASS" (i.e. the shared "FAIL" string is
# renamed to "PASS" in place).
$ws.Range("E7").Value = "SKIP"
$ws.Range("E12").Value = "PASS"

# --- view/selection edits --------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D20").Select() | Out-Null
